$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1
$ws.Range("B1").Value = "192.168.100.2455"
$ws.Range("D1").Value = "Nonedafdaf"

# Update row 2
$ws.Range("B2").Value = "192.168.100.244"
$ws.Range("D2").ClearContents()

# Remove row 5 entirely (shifts nothing below it, just deletes its contents/row)
$ws.Rows("5:5").Delete()
